# Insert a new weekly price-report row for Feria Lagunitas de Puerto Montt -
# Zanahoria. The new observation (2021-11-23) is inserted right after the
# existing row for 2021-05-26 (current row 169), pushing all subsequent rows
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at position 169 (shifts old rows 169..249 down to 170..250)
$ws.Rows.Item(169).Insert()

# Populate the newly inserted row with the new data point
$ws.Cells.Item(169, 1).Value = 4
$ws.Cells.Item(169, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(169, 3).Value = 'Los Lagos'
$ws.Cells.Item(169, 4).Value = 44523
$ws.Cells.Item(169, 5).Value = 10
$ws.Cells.Item(169, 6).Value = 100114013
$ws.Cells.Item(169, 7).Value = 'Zanahoria'
$ws.Cells.Item(169, 8).Value = 'Sin especificar'
$ws.Cells.Item(169, 9).Value = 'Primera'
$ws.Cells.Item(169, 10).Value = 750
$ws.Cells.Item(169, 11).Value = 11000
$ws.Cells.Item(169, 12).Value = 11000
$ws.Cells.Item(169, 13).Value = 11000
$ws.Cells.Item(169, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(169, 15).Value = 'Región de Ñuble'
$ws.Cells.Item(169, 16).Value = 550
$ws.Cells.Item(169, 17).Value = 20
$ws.Cells.Item(169, 18).Value = 'Hortaliza'
